$wb = $excel.ActiveWorkbook
$ws0 = $wb.Worksheets.Item(1)

# Insert a new worksheet ("Info") before the existing "{R-TITLE}" sheet,
# explaining the multi-sheet print-preview workaround.
$info = $wb.Worksheets.Add($ws0)
$info.Name = "Info"

$info.Range("A1").Value = "Multi sheet issues with print preview"
$info.Range("A1").Font.Bold = $true
$info.Range("A1").Font.Name = "Calibri"

$info.Range("A2").Value = "If you have problems with printing the document:"
$info.Range("A3").Value = "1)Create 1st sheet like this"
$info.Range("A4").Value = "2)or call macro like in example 070 to save file with a new name (so the entire file would be recreated)"

$info.Range("A5").Value = "https://docs.microsoft.com/en-us/office/vba/api/excel.workbook.saveas"
$link = $info.Hyperlinks.Add($info.Range("A5"), "https://docs.microsoft.com/en-us/office/vba/api/excel.workbook.saveas")

$info.Range("K1").Select() | Out-Null

Write-Output "done"
